$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Bot"
$ws.Range("B2").Value = "Salut! Cum te numești?"

$ws.Range("A3").Value = "User"
$ws.Range("B3").Value = "ma numesc Alexandru"

$ws.Range("A4").Value = "Bot"
$ws.Range("B4").Value = "Sunt un asistent virtual simplu."

$ws.Range("A5").Value = "User"
$ws.Range("B5").Value = "ce faci?"

$ws.Range("A6").Value = "Bot"
$ws.Range("B6").Value = "Sunt bine, mulțumesc că întrebi!"

$ws.Range("A7").Value = "User"
$ws.Range("B7").Value = "ce poti zice despre python?"

$ws.Range("B8").Value = "Sunt un asistent virtual simplu."
